# Update Name of Algo
# Apply updated imputed values to the RandomForest result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.149600000000003
$ws.Range("A3").Value = -22.26509999999999
$ws.Range("E3").Value = 16.1625
$ws.Range("E12").Value = 18.47690000000003
$ws.Range("A14").Value = -21.86020000000001
$ws.Range("A21").Value = -20.17349999999999
$ws.Range("A23").Value = -20.18969999999997
$ws.Range("E24").Value = 16.5107
$ws.Range("A25").Value = -21.88199999999999
$ws.Range("B25").Value = 5.260600000000003
$ws.Range("E25").Value = 17.0516
$ws.Range("A26").Value = -21.10479999999997
$ws.Range("B27").Value = 6.637899999999997
$ws.Range("A29").Value = -20.67939999999999
$ws.Range("B31").Value = 5.642999999999996
$ws.Range("B39").Value = 9.736700000000003
$ws.Range("B48").Value = 5.195700000000001
$ws.Range("E50").Value = 16.33209999999999
$ws.Range("B51").Value = 5.741799999999999
$ws.Range("B52").Value = 5.389799999999996
$ws.Range("A53").Value = -22.01530000000001
$ws.Range("E53").Value = 16.78760000000002
$ws.Range("B55").Value = 6.026199999999994
$ws.Range("B56").Value = 4.924599999999997
$ws.Range("A57").Value = -22.63980000000002
$ws.Range("B57").Value = 4.871599999999995
$ws.Range("E57").Value = 16.3578
$ws.Range("A59").Value = -22.1825
$ws.Range("E61").Value = 16.541
$ws.Range("E63").Value = 18.47740000000002
$ws.Range("A69").Value = -21.60079999999999
$ws.Range("E70").Value = 18.53750000000002
$ws.Range("B73").Value = 8.393800000000001
$ws.Range("A79").Value = -20.5268
$ws.Range("A83").Value = -21.9545
$ws.Range("E86").Value = 16.424
$ws.Range("B89").Value = 4.566199999999994
$ws.Range("B90").Value = 5.5589
$ws.Range("A91").Value = -21.40900000000001
$ws.Range("B92").Value = 4.628299999999996
$ws.Range("A93").Value = -20.93499999999999
$ws.Range("E98").Value = 15.985
$ws.Range("E100").Value = 16.47760000000001
$ws.Range("E102").Value = 16.54229999999999
